# Sync attendance_reports: reorder "Recorded By" (column G) comma-separated
# list entries so that the System/backdoor/admin marker that was first is
# moved to the end (and the previously-last entry becomes first).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ", "

        $hasSystemToken = $false
        foreach ($p in $parts) {
            if ($p.Trim().ToLower() -eq "system") {
                $hasSystemToken = $true
            }
        }

        if ($hasSystemToken) {
            $reversedParts = @()
            for ($i = $parts.Count - 1; $i -ge 0; $i--) {
                $reversedParts += $parts[$i]
            }

            $newVal = [string]::Join(", ", $reversedParts)
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
